$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tbl = $shape.Table

# Resize columns (widths in points; 1 pt = 12700 EMU)
$tbl.Columns.Item(1).Width = 48   # 609600 EMU
$tbl.Columns.Item(2).Width = 216  # 2743200 EMU
$tbl.Columns.Item(3).Width = 132  # 1676400 EMU
$tbl.Columns.Item(4).Width = 48   # 609600 EMU
$tbl.Columns.Item(5).Width = 36   # 457200 EMU
$tbl.Columns.Item(6).Width = 48   # 609600 EMU

# Collapse every row height to 0 (auto-fit)
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $tbl.Rows.Item($i).Height = 0
}

# Force the frame's vertical extent to 0 (matches target cy="0")
$shape.Height = 0

# Update the text of the cell that held "7249141700003" (row 3, col 2)
$cell = $tbl.Cell(3, 2)
$cell.Shape.TextFrame.TextRange.Text = "sdfsdfsfdsdfsdfsdf`r`nsfsdfsdfsdfsd"
